# Liverpool_stats.xlsx update:
#  1. Rename the per-stat worksheet tabs to human-friendly, spaced names.
#  2. Bump every player's "Age" cell (column E, format "YY-DDD") forward by
#     one day across all the per-stat sheets (rows below the two header
#     rows, excluding the "Squad Total" / "Opponent Total" summary rows).

$wb = $excel.ActiveWorkbook

# --- 1. Rename worksheet tabs -------------------------------------------
$renames = @{
    "StandardStats"    = "Standard Stats"
    "ShootingStats"     = "Shooting Stats"
    "PassingStats"      = "Passing Stats"
    "PassTypes"         = "Pass Types"
    "GoalShotCreation"  = "Goal & Shot Creation"
    "DefensiveActions"  = "Defensive Actions"
    "PlayingTime"       = "Playing Time"
    "MiscStats"         = "Miscellaneous Stats"
}

foreach ($oldName in $renames.Keys) {
    $sheet = $wb.Worksheets.Item($oldName)
    $sheet.Name = $renames[$oldName]
}

# --- 2. Increment the "Age" day count on every stats sheet --------------
# Column E holds values like "33-255" (years-days). Every player row on
# every stats sheet (sheet index 2 through 10) moves one day forward.
$ageRegex = '^(\d+)-(\d+)$'

for ($i = 2; $i -le 10; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $lastRow = $ws.UsedRange.Rows.Count

    for ($r = 4; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 5)
        $val = $cell.Value()

        if ($val -match $ageRegex) {
            $years = [int]$matches[1]
            $days = [int]$matches[2] + 1
            $cell.Value = "{0}-{1:D3}" -f $years, $days
        }
    }
}
